$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 336151.66
$ws.Range("J17").Value = 336151.66
$ws.Range("L17").Value = 1008454.98
$ws.Range("N17").Value = -1008790.98
$ws.Range("H100").Value = 7724.44
$ws.Range("I100").Value = 763.1111
$ws.Range("K100").Value = 763.1111
$ws.Range("M100").Value = -222.1111
$ws.Range("H111").Value = 38513.895
$ws.Range("I111").Value = 14716.857
$ws.Range("J111").Value = 105145.6
$ws.Range("K111").Value = 44150.571
$ws.Range("L111").Value = 315436.8
$ws.Range("M111").Value = -41083.571
$ws.Range("N111").Value = -321570.8
$ws.Range("H114").Value = 97861
$ws.Range("J114").Value = 97861
$ws.Range("L114").Value = 97861
$ws.Range("N114").Value = -106539
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 14710825
$ws.Range("I132").Value = 16669194
$ws.Range("K132").Value = 50007582
$ws.Range("M132").Value = -50005052
$ws.Range("H135").Value = 4638.385
$ws.Range("I135").Value = 2860.8235
$ws.Range("K135").Value = 25747.4115
$ws.Range("M135").Value = -23212.4115
$ws.Range("H137").Value = 6693.25
$ws.Range("I137").Value = 2617.8823
$ws.Range("K137").Value = 7853.646900000001
$ws.Range("M137").Value = -5303.646900000001
$ws.Range("H138").Value = 1224805.5
$ws.Range("I138").Value = 501000
$ws.Range("K138").Value = 1503000
$ws.Range("M138").Value = -1497860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8165.8823
$ws.Range("I32").Value = 7988.516
$ws.Range("K32").Value = 7988.516
$ws.Range("M32").Value = -7701.516
$ws.Range("H45").Value = 2245
$ws.Range("I45").Value = 2054.5833
$ws.Range("J45").Value = 2571.4285
$ws.Range("K45").Value = 2054.5833
$ws.Range("L45").Value = 2571.4285
$ws.Range("M45").Value = -1677.5833
$ws.Range("N45").Value = -3325.4285
$ws.Range("H61").Value = 3397.4285
$ws.Range("I61").Value = 1520.725
$ws.Range("J61").Value = 8089.1875
$ws.Range("K61").Value = 1520.725
$ws.Range("L61").Value = 8089.1875
$ws.Range("M61").Value = -1308.725
$ws.Range("N61").Value = -8513.1875
$ws.Range("H74").Value = 2129.9524
$ws.Range("I74").Value = 1593.5
$ws.Range("J74").Value = 3846.6
$ws.Range("K74").Value = 1593.5
$ws.Range("L74").Value = 3846.6
$ws.Range("M74").Value = -719.5
$ws.Range("N74").Value = -5594.6
$ws.Range("H77").Value = 2129.9524
$ws.Range("I77").Value = 1593.5
$ws.Range("J77").Value = 3846.6
$ws.Range("K77").Value = 7967.5
$ws.Range("L77").Value = 19233
$ws.Range("M77").Value = -3599.5
$ws.Range("N77").Value = -27969
$ws.Range("H106").Value = 89998
$ws.Range("J106").Value = 89998
$ws.Range("L106").Value = 89998
$ws.Range("N106").Value = -92522
$ws.Range("H132").Value = 3914.6128
$ws.Range("I132").Value = 1424.5186
$ws.Range("K132").Value = 4273.5558
$ws.Range("M132").Value = -1743.5558
$ws.Range("H136").Value = 3397.4285
$ws.Range("I136").Value = 1520.725
$ws.Range("J136").Value = 8089.1875
$ws.Range("K136").Value = 4562.174999999999
$ws.Range("L136").Value = 24267.5625
$ws.Range("M136").Value = -2012.174999999999
$ws.Range("N136").Value = -29367.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2498.16
$ws.Range("I20").Value = 2212.2
$ws.Range("J20").Value = 2927.1
$ws.Range("K20").Value = 2212.2
$ws.Range("L20").Value = 2927.1
$ws.Range("M20").Value = -1965.2
$ws.Range("N20").Value = -3421.1
$ws.Range("H105").Value = 2835.4688
$ws.Range("I105").Value = 2109.5925
$ws.Range("K105").Value = 2109.5925
$ws.Range("M105").Value = -362.5925000000002
$ws.Range("H134").Value = 4559.8555
$ws.Range("I134").Value = 1344.6279
$ws.Range("J134").Value = 8749.394
$ws.Range("K134").Value = 4033.8837
$ws.Range("L134").Value = 26248.182
$ws.Range("M134").Value = -1498.8837
$ws.Range("N134").Value = -31318.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2339.0833
$ws.Range("I58").Value = 2214.6924
$ws.Range("J58").Value = 2486.0908
$ws.Range("K58").Value = 2214.6924
$ws.Range("L58").Value = 2486.0908
$ws.Range("M58").Value = -2011.6924
$ws.Range("N58").Value = -2892.0908
$ws.Range("H94").Value = 1534.4667
$ws.Range("I94").Value = 1063
$ws.Range("K94").Value = 1063
$ws.Range("M94").Value = -612
$ws.Range("H112").Value = 100702
$ws.Range("J112").Value = 100702
$ws.Range("L112").Value = 100702
$ws.Range("N112").Value = -103656
$ws.Range("H132").Value = 1382672.2
$ws.Range("I132").Value = 1670179.2
$ws.Range("J132").Value = 2638.6
$ws.Range("K132").Value = 5010537.6
$ws.Range("L132").Value = 7915.799999999999
$ws.Range("M132").Value = -5008007.6
$ws.Range("N132").Value = -12975.8
$ws.Range("H134").Value = 3240.8333
$ws.Range("I134").Value = 1846.591
$ws.Range("J134").Value = 7075
$ws.Range("K134").Value = 5539.772999999999
$ws.Range("L134").Value = 21225
$ws.Range("M134").Value = -3004.772999999999
$ws.Range("N134").Value = -26295
$ws.Range("H136").Value = 2339.0833
$ws.Range("I136").Value = 2214.6924
$ws.Range("J136").Value = 2486.0908
$ws.Range("K136").Value = 6644.0772
$ws.Range("L136").Value = 7458.2724
$ws.Range("M136").Value = -4094.0772
$ws.Range("N136").Value = -12558.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7728.6
$ws.Range("I5").Value = 337.375
$ws.Range("K5").Value = 1012.125
$ws.Range("M5").Value = -900.125
$ws.Range("H80").Value = 11749
$ws.Range("J80").Value = 12498.5
$ws.Range("L80").Value = 37495.5
$ws.Range("N80").Value = -39367.5
$ws.Range("H83").Value = 11749
$ws.Range("J83").Value = 12498.5
$ws.Range("L83").Value = 112486.5
$ws.Range("N83").Value = -121846.5
$ws.Range("H114").Value = 3388.5
$ws.Range("J114").Value = 4351.3335
$ws.Range("L114").Value = 13054.0005
$ws.Range("N114").Value = -19562.0005
$ws.Range("H129").Value = 1857
$ws.Range("J129").Value = 3661.1667
$ws.Range("L129").Value = 10983.5001
$ws.Range("N129").Value = -20983.5001
$ws.Range("H132").Value = 3284.6296
$ws.Range("I132").Value = 1225
$ws.Range("J132").Value = 3642.8262
$ws.Range("K132").Value = 11025
$ws.Range("L132").Value = 32785.4358
$ws.Range("M132").Value = -8495
$ws.Range("N132").Value = -37845.4358
$ws.Range("H135").Value = 7728.6
$ws.Range("I135").Value = 337.375
$ws.Range("K135").Value = 3036.375
$ws.Range("M135").Value = -501.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 61429.277
$ws.Range("I43").Value = 30147.6
$ws.Range("J43").Value = 73460.69500000001
$ws.Range("K43").Value = 30147.6
$ws.Range("L43").Value = 73460.69500000001
$ws.Range("M43").Value = -29996.6
$ws.Range("N43").Value = -73762.69500000001
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2078.875
$ws.Range("I22").Value = 639.5
$ws.Range("J22").Value = 3107
$ws.Range("K22").Value = 639.5
$ws.Range("L22").Value = 3107
$ws.Range("M22").Value = -344.5
$ws.Range("N22").Value = -3697
$ws.Range("H27").Value = 2078.875
$ws.Range("I27").Value = 639.5
$ws.Range("J27").Value = 3107
$ws.Range("K27").Value = 639.5
$ws.Range("L27").Value = 3107
$ws.Range("M27").Value = -532.5
$ws.Range("N27").Value = -3321
$ws.Range("H132").Value = 2617.1404
$ws.Range("I132").Value = 2550.8113
$ws.Range("J132").Value = 3496
$ws.Range("K132").Value = 7652.4339
$ws.Range("L132").Value = 10488
$ws.Range("M132").Value = -5122.4339
$ws.Range("N132").Value = -15548
$ws.Range("H136").Value = 3530.077
$ws.Range("I136").Value = 3059.0322
$ws.Range("K136").Value = 9177.096600000001
$ws.Range("M136").Value = -6627.096600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 56783.4
$ws.Range("I81").Value = 147435.86
$ws.Range("K81").Value = 294871.72
$ws.Range("M81").Value = -293810.72
$ws.Range("H84").Value = 56783.4
$ws.Range("I84").Value = 147435.86
$ws.Range("K84").Value = 1474358.6
$ws.Range("M84").Value = -1469054.6
$ws.Range("H122").Value = 3754.6667
$ws.Range("I122").Value = 3538.4092
$ws.Range("J122").Value = 4349.375
$ws.Range("K122").Value = 10615.2276
$ws.Range("L122").Value = 13048.125
$ws.Range("M122").Value = -8165.2276
$ws.Range("N122").Value = -17948.125
$ws.Range("H126").Value = 2405.6775
$ws.Range("I126").Value = 2110.739
$ws.Range("K126").Value = 6332.217000000001
$ws.Range("M126").Value = -3862.217000000001
$ws.Range("H132").Value = 3353.2666
$ws.Range("I132").Value = 3450.7856
$ws.Range("J132").Value = 1988
$ws.Range("K132").Value = 10352.3568
$ws.Range("L132").Value = 5964
$ws.Range("M132").Value = -7822.356800000001
$ws.Range("N132").Value = -11024
$ws.Range("H136").Value = 8023.483
$ws.Range("I136").Value = 9026.48
$ws.Range("K136").Value = 27079.44
$ws.Range("M136").Value = -24529.44
